# Apply weekly stock-history update to ADANIENSOL.NS sheet:
#  - row 44: O44 (isPivot) 0 -> 3
#  - rows 46/47: R (backup) blank -> 0 (now "completed"/confirmed)
#  - append 24 new weekly rows (48-71) of OHLCV + date-part data,
#    whose own R (backup) column stays blank until later confirmed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing-row fixups -------------------------------------------------
$ws.Cells.Item(44, 15).Value = 3      # O44: isPivot 0 -> 3
$ws.Cells.Item(46, 18).Value = 0      # R46: backup blank -> 0
$ws.Cells.Item(47, 18).Value = 0      # R47: backup blank -> 0

# --- new weekly rows (A:Q); R (backup) intentionally left blank ----------
$newRows = @(
    @(45474, 1003.650024414062, 1063.349975585938, 996.0499877929688, 1012.650024414062, 1012.650024414062, 5665904, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
    @(45481, 1010, 1033.199951171875, 985.4000244140625, 1004.799987792969, 1004.799987792969, 3159643, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(45488, 1006.5, 1047.949951171875, 998.5, 1010.700012207031, 1010.700012207031, 3631455, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 0),
    @(45495, 1004.950012207031, 1149.699951171875, 975, 1051.949951171875, 1051.949951171875, 12467766, 2024, 7, 22, 0, 0, 0, 30, 2, 0, 0),
    @(45502, 1059.150024414062, 1348, 1047.849975585938, 1261.800048828125, 1261.800048828125, 58300860, 2024, 7, 29, 0, 0, 0, 31, 1, 0, 0),
    @(45509, 1200, 1231.949951171875, 1094.300048828125, 1103.800048828125, 1103.800048828125, 27769921, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(45516, 1072, 1138, 1036.050048828125, 1085.300048828125, 1085.300048828125, 19210757, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(45523, 1089.800048828125, 1115.5, 1060, 1065.650024414062, 1065.650024414062, 8338037, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(45530, 1066.550048828125, 1076, 984, 1007.400024414062, 1007.400024414062, 12181819, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(45537, 1017, 1049, 983, 985.3499755859375, 985.3499755859375, 8068349, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 1),
    @(45544, 985.2999877929688, 1020, 968.2999877929688, 983.2999877929688, 983.2999877929688, 6951647, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 2),
    @(45551, 990.5, 1019.900024414062, 967.0499877929688, 1007.150024414062, 1007.150024414062, 14159284, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
    @(45558, 1025.199951171875, 1068, 1001.049987792969, 1010.200012207031, 1010.200012207031, 9725871, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0),
    @(45565, 1009.900024414062, 1053.449951171875, 955, 962.75, 962.75, 11388322, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0),
    @(45572, 976, 1017.450012207031, 915, 976.7999877929688, 976.7999877929688, 14093009, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0),
    @(45579, 994, 1046, 974.0499877929688, 1043.300048828125, 1043.300048828125, 10532710, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0),
    @(45586, 1043.300048828125, 1048.699951171875, 905.0999755859375, 919.3499755859375, 919.3499755859375, 11448149, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0),
    @(45593, 924.9500122070312, 990, 891.0499877929688, 978.75, 978.75, 5820533, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0),
    @(45600, 978.9500122070312, 1090.949951171875, 934, 936.25, 936.25, 18677702, 2024, 11, 4, 0, 0, 0, 45, 1, 0, 0),
    @(45607, 936.25, 937.9500122070312, 875.8499755859375, 879.7000122070312, 879.7000122070312, 11799539, 2024, 11, 11, 0, 0, 0, 46, 0, 0, 0),
    @(45614, 888, 889, 627.5, 648.9500122070312, 648.9500122070312, 47598159, 2024, 11, 18, 0, 0, 0, 47, 0, 0, 0),
    @(45621, 678, 869.75, 588, 840.5, 840.5, 70074158, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0),
    @(45628, 844.9000244140625, 855.8499755859375, 797, 808.3499755859375, 808.3499755859375, 21713159, 2024, 12, 2, 0, 0, 0, 49, 0, 0, 0),
    @(45635, 808.3499755859375, 837.5, 773, 832.4000244140625, 832.4000244140625, 16446055, 2024, 12, 9, 0, 0, 0, 50, 0, 0, 0)
)

$startRow = 48
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]          # A: Datetime (serial)
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $vals[1]          # B: Open
    $ws.Cells.Item($r, 3).Value = $vals[2]          # C: High
    $ws.Cells.Item($r, 4).Value = $vals[3]          # D: Low
    $ws.Cells.Item($r, 5).Value = $vals[4]          # E: Close
    $ws.Cells.Item($r, 6).Value = $vals[5]          # F: Adj Close
    $ws.Cells.Item($r, 7).Value = $vals[6]          # G: Volume
    $ws.Cells.Item($r, 8).Value = $vals[7]          # H: Year
    $ws.Cells.Item($r, 9).Value = $vals[8]          # I: Month
    $ws.Cells.Item($r, 10).Value = $vals[9]         # J: Day
    $ws.Cells.Item($r, 11).Value = $vals[10]        # K: Hour
    $ws.Cells.Item($r, 12).Value = $vals[11]        # L: Minute
    $ws.Cells.Item($r, 13).Value = $vals[12]        # M: Second
    $ws.Cells.Item($r, 14).Value = $vals[13]        # N: Week
    $ws.Cells.Item($r, 15).Value = $vals[14]        # O: isPivot
    $ws.Cells.Item($r, 16).Value = $vals[15]        # P: two_line_structure
    $ws.Cells.Item($r, 17).Value = $vals[16]        # Q: detect_structure
    # R (backup) left unset -> blank, matching the not-yet-confirmed rows
}

Write-Host "Applied ADANIENSOL.NS weekly-data update: O44, R46, R47, rows 48-71"
